# Adds a block of new paragraphs (a greeting / training-status request
# email body) right after the existing table and before the document's
# final (bookmark) paragraph.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Common pPr used by all the "NormalWeb" styled paragraphs being added.
$pPrCommon = @"
<w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>
"@

$rPrCommon = @"
<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>
"@

# Paragraph fragments, in the order they must appear.
$paraXmls = @()

# 1) A bare empty paragraph (no formatting at all).
$paraXmls += "<w:p $wNs/>"

# 2) "Hi Manoj,Kautak,Fathima,"
$paraXmls += "<w:p $wNs>$pPrCommon" +
  "<w:r>$rPrCommon<w:t xml:space=`"preserve`">Hi </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r>$rPrCommon<w:t>Manoj</w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r>$rPrCommon<w:t>,Kautak,Fathima</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r>$rPrCommon<w:t>,</w:t></w:r>" +
  "</w:p>"

# 3) A single space.
$paraXmls += "<w:p $wNs>$pPrCommon<w:r>$rPrCommon<w:t> </w:t></w:r></w:p>"

# 4) Training details request.
$paraXmls += "<w:p $wNs>$pPrCommon<w:r>$rPrCommon<w:t>Pleas share the details of the trainings completed so far and Knowledge gained so far post training with tech stack and ranking on scale of 10.</w:t></w:r></w:p>"

# 5) Duration sentence.
$paraXmls += "<w:p $wNs>$pPrCommon<w:r>$rPrCommon<w:t>Include duration of training as well.</w:t></w:r></w:p>"

# 6) A single space.
$paraXmls += "<w:p $wNs>$pPrCommon<w:r>$rPrCommon<w:t> </w:t></w:r></w:p>"

# 7) Closing sentence.
$paraXmls += "<w:p $wNs>$pPrCommon<w:r>$rPrCommon<w:t>This Info need to be shared with customer and during customer evaluation questions will be focused on same.</w:t></w:r></w:p>"

$count = $paraXmls.Count

# Locate the document's current last paragraph (the one holding the
# "_GoBack" bookmark and the trailing line break) and insert the new
# paragraphs immediately before it, preserving their order.
$n0 = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n0)
$insertionRange = $lastPara.Range
$insertionRange.Collapse(1)

for ($i = 0; $i -lt $count; $i++) {
    $insertionRange.InsertParagraphBefore()
}

# The freshly minted (still empty) paragraphs now sit at indices
# n0 .. n0+count-1; fill each one in with its real content/formatting.
for ($i = 0; $i -lt $count; $i++) {
    $p = $d.Paragraphs.Item($n0 + $i)
    $p.Range.InsertXML($paraXmls[$i])
}
